$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 3901
$ws.Range("J38").Value = 9800
$ws.Range("L38").Value = 29400
$ws.Range("N38").Value = -30144

# Row 57
$ws.Range("H57").Value = 35208.25
$ws.Range("J57").Value = 35779
$ws.Range("L57").Value = 107337
$ws.Range("N57").Value = -108335

# Row 86
$ws.Range("H86").Value = 5186.75
$ws.Range("I86").Value = 5356.4614
$ws.Range("J86").Value = 4871.5713
$ws.Range("K86").Value = 5356.4614
$ws.Range("L86").Value = 4871.5713
$ws.Range("M86").Value = -4233.4614
$ws.Range("N86").Value = -7117.5713

# Row 89
$ws.Range("H89").Value = 5186.75
$ws.Range("I89").Value = 5356.4614
$ws.Range("J89").Value = 4871.5713
$ws.Range("K89").Value = 26782.307
$ws.Range("L89").Value = 24357.8565
$ws.Range("M89").Value = -21166.307
$ws.Range("N89").Value = -35589.85649999999

# Row 132
$ws.Range("H132").Value = 1874.6825
$ws.Range("I132").Value = 1923.2069
$ws.Range("K132").Value = 5769.620699999999
$ws.Range("M132").Value = -3239.620699999999

# Row 137
$ws.Range("H137").Value = 4489.1665
$ws.Range("I137").Value = 4489.1665
$ws.Range("K137").Value = 13467.4995
$ws.Range("M137").Value = -10917.4995

# Row 138
$ws.Range("H138").Value = 9595.59
$ws.Range("I138").Value = 7657.375
$ws.Range("J138").Value = 9764.130999999999
$ws.Range("K138").Value = 22972.125
$ws.Range("L138").Value = 29292.393
$ws.Range("M138").Value = -17832.125
$ws.Range("N138").Value = -39572.393

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18322.057
$ws.Range("I32").Value = 17962.885
$ws.Range("J32").Value = 36999
$ws.Range("K32").Value = 17962.885
$ws.Range("L32").Value = 36999
$ws.Range("M32").Value = -17675.885
$ws.Range("N32").Value = -37573

# Row 45
$ws.Range("H45").Value = 3206
$ws.Range("I45").Value = 2391.2144
$ws.Range("J45").Value = 4835.5713
$ws.Range("K45").Value = 2391.2144
$ws.Range("L45").Value = 4835.5713
$ws.Range("M45").Value = -2014.2144
$ws.Range("N45").Value = -5589.5713

# Row 74
$ws.Range("H74").Value = 2886.3333
$ws.Range("I74").Value = 2853.4119
$ws.Range("K74").Value = 2853.4119
$ws.Range("M74").Value = -1979.4119

# Row 77
$ws.Range("H77").Value = 2886.3333
$ws.Range("I77").Value = 2853.4119
$ws.Range("K77").Value = 14267.0595
$ws.Range("M77").Value = -9899.059499999999

# Row 102
$ws.Range("H102").Value = 50004332
$ws.Range("I102").Value = 71432020
$ws.Range("K102").Value = 71432020
$ws.Range("M102").Value = -71430398

# Row 122
$ws.Range("H122").Value = 6038.6665
$ws.Range("I122").Value = 6038.6665
$ws.Range("K122").Value = 18115.9995
$ws.Range("M122").Value = -15665.9995

# Row 132
$ws.Range("H132").Value = 14289078
$ws.Range("I132").Value = 3924.3333
$ws.Range("K132").Value = 11772.9999
$ws.Range("M132").Value = -9242.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3258.1277
$ws.Range("I20").Value = 3006.3462
$ws.Range("J20").Value = 3569.8572
$ws.Range("K20").Value = 3006.3462
$ws.Range("L20").Value = 3569.8572
$ws.Range("M20").Value = -2759.3462
$ws.Range("N20").Value = -4063.8572

# Row 57
$ws.Range("H57").Value = 109902.38
$ws.Range("J57").Value = 109902.38
$ws.Range("L57").Value = 109902.38
$ws.Range("N57").Value = -111342.38

# Row 60
$ws.Range("H60").Value = 188997.67
$ws.Range("J60").Value = 188997.67
$ws.Range("L60").Value = 188997.67
$ws.Range("N60").Value = -190195.67

# Row 132
$ws.Range("H132").Value = 120739.43
$ws.Range("J132").Value = 120739.43
$ws.Range("L132").Value = 120739.43
$ws.Range("N132").Value = -130859.43

# Row 134
$ws.Range("H134").Value = 12502424
$ws.Range("I134").Value = 2770.7144
$ws.Range("J134").Value = 100000000
$ws.Range("K134").Value = 8312.143199999999
$ws.Range("L134").Value = 300000000
$ws.Range("M134").Value = -5777.143199999999
$ws.Range("N134").Value = -300005070

# Row 136
$ws.Range("H136").Value = 109902.38
$ws.Range("J136").Value = 109902.38
$ws.Range("L136").Value = 109902.38
$ws.Range("N136").Value = -120102.38

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1579.6923
$ws.Range("J22").Value = 2433.1667
$ws.Range("L22").Value = 2433.1667
$ws.Range("N22").Value = -3133.1667

# Row 31
$ws.Range("H31").Value = 35718756
$ws.Range("I31").Value = 55558676
$ws.Range("K31").Value = 55558676
$ws.Range("M31").Value = -55558381

# Row 34
$ws.Range("H34").Value = 35718756
$ws.Range("I34").Value = 55558676
$ws.Range("K34").Value = 55558676
$ws.Range("M34").Value = -55558474

# Row 52
$ws.Range("H52").Value = 68063.86
$ws.Range("I52").Value = 30000
$ws.Range("J52").Value = 74407.836
$ws.Range("K52").Value = 30000
$ws.Range("L52").Value = 74407.836
$ws.Range("M52").Value = -29706
$ws.Range("N52").Value = -74995.836

# Row 132
$ws.Range("H132").Value = 2895.0454
$ws.Range("I132").Value = 2689.55
$ws.Range("J132").Value = 4950
$ws.Range("K132").Value = 8068.650000000001
$ws.Range("L132").Value = 14850
$ws.Range("M132").Value = -5538.650000000001
$ws.Range("N132").Value = -19910

# Row 135
$ws.Range("H135").Value = 110498.086
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 110498.086
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 110498.086
$ws.Range("M135").Value = ""
$ws.Range("N135").Value = -120638.086

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 128969350
$ws.Range("I4").Value = 145964980
$ws.Range("K4").Value = 437894940
$ws.Range("M4").Value = -437894828

# Row 121
$ws.Range("H121").Value = 2456
$ws.Range("I121").Value = 230.44444
$ws.Range("J121").Value = 3886.7144
$ws.Range("K121").Value = 691.33332
$ws.Range("L121").Value = 11660.1432
$ws.Range("M121").Value = 618.66668
$ws.Range("N121").Value = -14280.1432

# Row 122
$ws.Range("H122").Value = 84093.5
$ws.Range("J122").Value = 2124.6667
$ws.Range("L122").Value = 19122.0003
$ws.Range("N122").Value = -24022.0003

# Row 139
$ws.Range("H139").Value = 6035.75
$ws.Range("I139").Value = 2761.2222
$ws.Range("J139").Value = 10245.857
$ws.Range("K139").Value = 8283.6666
$ws.Range("L139").Value = 30737.571
$ws.Range("M139").Value = -3143.6666
$ws.Range("N139").Value = -41017.571

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8007.8
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = ""

# Row 73
$ws.Range("H73").Value = 8007.8
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = ""

# Row 132
$ws.Range("H132").Value = 8305185
$ws.Range("I132").Value = 5738.6
$ws.Range("K132").Value = 17215.8
$ws.Range("M132").Value = -14685.8

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 16707466
$ws.Range("I100").Value = 5198.375
$ws.Range("J100").Value = 35795772
$ws.Range("K100").Value = 5198.375
$ws.Range("L100").Value = 35795772
$ws.Range("M100").Value = -4657.375
$ws.Range("N100").Value = -35796854

# Row 136
$ws.Range("H136").Value = 10365.583
$ws.Range("I136").Value = 13782.571
$ws.Range("J136").Value = 5581.8
$ws.Range("K136").Value = 41347.713
$ws.Range("L136").Value = 16745.4
$ws.Range("M136").Value = -38797.713
$ws.Range("N136").Value = -21845.4

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = ""

# Row 100
$ws.Range("H100").Value = 2355.375
$ws.Range("I100").Value = 1627
$ws.Range("J100").Value = 3291.8572
$ws.Range("K100").Value = 3254
$ws.Range("L100").Value = 6583.7144
$ws.Range("M100").Value = -2713
$ws.Range("N100").Value = -7665.7144

# Row 132
$ws.Range("H132").Value = 1670497
$ws.Range("I132").Value = 4596.4
$ws.Range("K132").Value = 13789.2
$ws.Range("M132").Value = -11259.2

# Row 136
$ws.Range("H136").Value = 406914.53
$ws.Range("I136").Value = 7202.625
$ws.Range("K136").Value = 21607.625
$ws.Range("M136").Value = -19057.875
